$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_5_3_24"
$ws.Range("B2").Value = 0.5550301287191901
$ws.Range("C2").Value = 0.006776442350362566
$ws.Range("D2").Value = -2.671469480003284
$ws.Range("E2").Value = -1.368327815690432
$ws.Range("F2").Value = 0.4924505054950714
$ws.Range("G2").Value = 0.4854681193828583
$ws.Range("H2").Value = 2.482096910476685
$ws.Range("I2").Value = 1.425058007240295

$ws.Range("A3").Value = "model_5_3_12"
$ws.Range("B3").Value = 0.5579973425867283
$ws.Range("C3").Value = -0.8357295935230498
$ws.Range("D3").Value = -1.380530130829171
$ws.Range("E3").Value = -1.048094870723453
$ws.Range("F3").Value = 0.4891667068004608
$ws.Range("G3").Value = 0.8972684144973755
$ws.Range("H3").Value = 1.609357237815857
$ws.Range("I3").Value = 1.232369065284729

$ws.Range("A4").Value = "model_5_3_23"
$ws.Range("B4").Value = 0.5580358416448823
$ws.Range("C4").Value = 0.01213577285578915
$ws.Range("D4").Value = -2.612549710523224
$ws.Range("E4").Value = -1.334870205831018
$ws.Range("F4").Value = 0.4891240894794464
$ws.Range("G4").Value = 0.4828485250473022
$ws.Range("H4").Value = 2.442264080047607
$ws.Range("I4").Value = 1.404926180839539

$ws.Range("A5").Value = "model_5_3_21"
$ws.Range("B5").Value = 0.5621204493605951
$ws.Range("C5").Value = -0.4699587757138119
$ws.Range("D5").Value = -2.090886521625336
$ws.Range("E5").Value = -1.266378158539548
$ws.Range("F5").Value = 0.4846036434173584
$ws.Range("G5").Value = 0.7184867858886719
$ws.Range("H5").Value = 2.089593648910522
$ws.Range("I5").Value = 1.363713383674622

$ws.Range("A6").Value = "model_5_3_22"
$ws.Range("B6").Value = 0.5624276081260009
$ws.Range("C6").Value = -0.01351695435847211
$ws.Range("D6").Value = -2.520628980037091
$ws.Range("E6").Value = -1.297301758149455
$ws.Range("F6").Value = 0.484263688325882
$ws.Range("G6").Value = 0.4953870177268982
$ws.Range("H6").Value = 2.380120992660522
$ws.Range("I6").Value = 1.382320642471313

$ws.Range("A7").Value = "model_5_3_19"
$ws.Range("B7").Value = 0.5648157283188475
$ws.Range("C7").Value = -0.7718050275667088
$ws.Range("D7").Value = -1.631313261733653
$ws.Range("E7").Value = -1.153199446267144
$ws.Range("F7").Value = 0.4816207885742188
$ws.Range("G7").Value = 0.8660233020782471
$ws.Range("H7").Value = 1.778899312019348
$ws.Range("I7").Value = 1.295612096786499

$ws.Range("A8").Value = "model_5_3_13"
$ws.Range("B8").Value = 0.5684172925324882
$ws.Range("C8").Value = -0.7471151924975032
$ws.Range("D8").Value = -1.382223713733608
$ws.Range("E8").Value = -1.010881796376387
$ws.Range("F8").Value = 0.4776349067687988
$ws.Range("G8").Value = 0.8539553880691528
$ws.Range("H8").Value = 1.610502362251282
$ws.Range("I8").Value = 1.209977388381958

$ws.Range("A9").Value = "model_5_3_14"
$ws.Range("B9").Value = 0.5711367082256725
$ws.Range("C9").Value = -0.7346159519900224
$ws.Range("D9").Value = -1.403690777567814
$ws.Range("E9").Value = -1.01685698032231
$ws.Range("F9").Value = 0.474625289440155
$ws.Range("G9").Value = 0.8478460311889648
$ws.Range("H9").Value = 1.625015139579773
$ws.Range("I9").Value = 1.21357274055481

$ws.Range("A10").Value = "model_5_3_18"
$ws.Range("B10").Value = 0.5716242395529803
$ws.Range("C10").Value = -0.71219527415607
$ws.Range("D10").Value = -1.569871242808115
$ws.Range("E10").Value = -1.095078087640012
$ws.Range("F10").Value = 0.4740857481956482
$ws.Range("G10").Value = 0.8368872404098511
$ws.Range("H10").Value = 1.737361431121826
$ws.Range("I10").Value = 1.260639548301697

$ws.Range("A11").Value = "model_5_3_15"
$ws.Range("B11").Value = 0.5724229241434513
$ws.Range("C11").Value = -0.7138301189459089
$ws.Range("D11").Value = -1.446872000362969
$ws.Range("E11").Value = -1.03074892579841
$ws.Range("F11").Value = 0.4732018411159515
$ws.Range("G11").Value = 0.83768630027771
$ws.Range("H11").Value = 1.654207706451416
$ws.Range("I11").Value = 1.2219318151474

$ws.Range("A12").Value = "model_5_3_17"
$ws.Range("B12").Value = 0.575468971043079
$ws.Range("C12").Value = -0.6718694855355634
$ws.Range("D12").Value = -1.525335147593555
$ws.Range("E12").Value = -1.054189394539035
$ws.Range("F12").Value = 0.4698307514190674
$ws.Range("G12").Value = 0.8171766996383667
$ws.Range("H12").Value = 1.707252740859985
$ws.Range("I12").Value = 1.23603630065918

$ws.Range("A13").Value = "model_5_3_16"
$ws.Range("B13").Value = 0.5758054137512654
$ws.Range("C13").Value = -0.6570649503200447
$ws.Range("D13").Value = -1.494091828594903
$ws.Range("E13").Value = -1.031303357081398
$ws.Range("F13").Value = 0.4694584310054779
$ws.Range("G13").Value = 0.8099405765533447
$ws.Range("H13").Value = 1.68613076210022
$ws.Range("I13").Value = 1.222265481948853

$ws.Range("A14").Value = "model_5_3_11"
$ws.Range("B14").Value = 0.5766610666220124
$ws.Range("C14").Value = -0.664623388881517
$ws.Range("D14").Value = -1.26688681594009
$ws.Range("E14").Value = -0.9144251740897891
$ws.Range("F14").Value = 0.4685114324092865
$ws.Range("G14").Value = 0.8136350512504578
$ws.Range("H14").Value = 1.532528638839722
$ws.Range("I14").Value = 1.151938080787659

$ws.Range("A15").Value = "model_5_3_20"
$ws.Range("B15").Value = 0.5833880741212497
$ws.Range("C15").Value = -0.3054468308517095
$ws.Range("D15").Value = -1.960538699105742
$ws.Range("E15").Value = -1.126712536339167
$ws.Range("F15").Value = 0.4610666632652283
$ws.Range("G15").Value = 0.6380766034126282
$ws.Range("H15").Value = 2.001471996307373
$ws.Range("I15").Value = 1.279674530029297

$ws.Range("A16").Value = "model_5_3_10"
$ws.Range("B16").Value = 0.5889869893647994
$ws.Range("C16").Value = -0.5845158621796052
$ws.Range("D16").Value = -1.137830138801714
$ws.Range("E16").Value = -0.8117393631244987
$ws.Range("F16").Value = 0.4548702836036682
$ws.Range("G16").Value = 0.7744800448417664
$ws.Range("H16").Value = 1.445280075073242
$ws.Range("I16").Value = 1.090150475502014

$ws.Range("A17").Value = "model_5_3_9"
$ws.Range("B17").Value = 0.5947903309082923
$ws.Range("C17").Value = -0.5894879643179232
$ws.Range("D17").Value = -1.017278267084292
$ws.Range("E17").Value = -0.7501393564201948
$ws.Range("F17").Value = 0.4484476745128632
$ws.Range("G17").Value = 0.7769103050231934
$ws.Range("H17").Value = 1.363780975341797
$ws.Range("I17").Value = 1.053084850311279

$ws.Range("A18").Value = "model_5_3_8"
$ws.Range("B18").Value = 0.6026112360226544
$ws.Range("C18").Value = -0.5821645865941678
$ws.Range("D18").Value = -0.8911003367281183
$ws.Range("E18").Value = -0.6802765851878372
$ws.Range("F18").Value = 0.4397923052310944
$ws.Range("G18").Value = 0.7733308076858521
$ws.Range("H18").Value = 1.278478384017944
$ws.Range("I18").Value = 1.01104736328125

$ws.Range("A19").Value = "model_5_3_7"
$ws.Range("B19").Value = 0.613144983991541
$ws.Range("C19").Value = -0.5965954532019919
$ws.Range("D19").Value = -0.7051687723539604
$ws.Range("E19").Value = -0.588176495309725
$ws.Range("F19").Value = 0.4281345009803772
$ws.Range("G19").Value = 0.7803843021392822
$ws.Range("H19").Value = 1.152779340744019
$ws.Range("I19").Value = 0.9556294679641724

$ws.Range("A20").Value = "model_5_3_6"
$ws.Range("B20").Value = 0.6212918548933679
$ws.Range("C20").Value = -0.5558891093535625
$ws.Range("D20").Value = -0.5962758907312189
$ws.Range("E20").Value = -0.5130959610083703
$ws.Range("F20").Value = 0.419118344783783
$ws.Range("G20").Value = 0.7604878544807434
$ws.Range("H20").Value = 1.079162240028381
$ws.Range("I20").Value = 0.9104523062705994

$ws.Range("A21").Value = "model_5_3_5"
$ws.Range("B21").Value = 0.6297631725067594
$ws.Range("C21").Value = -0.5683919507471931
$ws.Range("D21").Value = -0.4302506504797676
$ws.Range("E21").Value = -0.4306914069280428
$ws.Range("F21").Value = 0.409743070602417
$ws.Range("G21").Value = 0.7665989398956299
$ws.Range("H21").Value = 0.9669208526611328
$ws.Range("I21").Value = 0.8608682751655579

$ws.Range("A22").Value = "model_5_3_4"
$ws.Range("B22").Value = 0.6351483482436855
$ws.Range("C22").Value = -0.5394192377716829
$ws.Range("D22").Value = -0.3326028794532103
$ws.Range("E22").Value = -0.3666026805125491
$ws.Range("F22").Value = 0.4037832617759705
$ws.Range("G22").Value = 0.7524377107620239
$ws.Range("H22").Value = 0.9009061455726624
$ws.Range("I22").Value = 0.8223051428794861

$ws.Range("A23").Value = "model_5_3_2"
$ws.Range("B23").Value = 0.6395712387529894
$ws.Range("C23").Value = -0.5794925361723628
$ws.Range("D23").Value = -0.1135658762688496
$ws.Range("E23").Value = -0.2680259874122277
$ws.Range("F23").Value = 0.3988884389400482
$ws.Range("G23").Value = 0.7720246911048889
$ws.Range("H23").Value = 0.7528261542320251
$ws.Range("I23").Value = 0.7629901170730591

$ws.Range("A24").Value = "model_5_3_3"
$ws.Range("B24").Value = 0.64077402888374
$ws.Range("C24").Value = -0.5494386517248586
$ws.Range("D24").Value = -0.2013364513512852
$ws.Range("E24").Value = -0.3015076396180443
$ws.Range("F24").Value = 0.397557258605957
$ws.Range("G24").Value = 0.7573350667953491
$ws.Range("H24").Value = 0.8121634721755981
$ws.Range("I24").Value = 0.7831364870071411

$ws.Range("A25").Value = "model_5_3_1"
$ws.Range("B25").Value = 0.6932146697000265
$ws.Range("C25").Value = 0.09680249171218636
$ws.Range("D25").Value = 0.4044617112107701
$ws.Range("E25").Value = 0.2967063077025471
$ws.Range("F25").Value = 0.3395209610462189
$ws.Range("G25").Value = 0.4414651095867157
$ws.Range("H25").Value = 0.402613639831543
$ws.Range("I25").Value = 0.4231822490692139

$ws.Range("A26").Value = "model_5_3_0"
$ws.Range("B26").Value = 0.7320359580327026
$ws.Range("C26").Value = 0.6547099464826277
$ws.Range("D26").Value = 0.5090547838950488
$ws.Range("E26").Value = 0.5919342933309111
$ws.Range("F26").Value = 0.2965572476387024
$ws.Range("G26").Value = 0.168770968914032
$ws.Range("H26").Value = 0.3319035172462463
$ws.Range("I26").Value = 0.2455392181873322
